$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension implicitly by writing data rows 2-21 (new rows 17-21 added)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Pros1"
$ws.Range("C2").Value = "Tyro3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 23.959728
$ws.Range("H2").Value = 71.879184
$ws.Range("I2").Value = 0.1471048789571275
$ws.Range("J2").Value = 0.1548696879699961
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.2063693333333333
$ws.Range("N2").Value = 0.619108
$ws.Range("O2").Value = 0.05249514260861875
$ws.Range("P2").Value = 0.05463357984752036
$ws.Range("Q2").Value = 4.944553094208
$ws.Range("R2").Value = 44.50097784787199
$ws.Range("S2").Value = 0.007722291599278008
$ws.Range("T2").Value = 0.008461085463669343

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Pros1"
$ws.Range("C3").Value = "Tyro3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 23.959728
$ws.Range("H3").Value = 71.879184
$ws.Range("I3").Value = 0.1471048789571275
$ws.Range("J3").Value = 0.1548696879699961
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.258868
$ws.Range("N3").Value = 9.776603999999999
$ws.Range("O3").Value = 0.8289736543672389
$ws.Range("P3").Value = 0.8627426479250582
$ws.Range("Q3").Value = 78.08159086790398
$ws.Range("R3").Value = 702.7343178111358
$ws.Range("S3").Value = 0.1219460690843403
$ws.Range("T3").Value = 0.1336126846825619

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Pros1"
$ws.Range("C4").Value = "Tyro3"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 23.959728
$ws.Range("H4").Value = 71.879184
$ws.Range("I4").Value = 0.1471048789571275
$ws.Range("J4").Value = 0.1548696879699961
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.004351333333333333
$ws.Range("N4").Value = 0.013054
$ws.Range("O4").Value = 0.001106869224130377
$ws.Range("P4").Value = 0.001151958545729551
$ws.Range("Q4").Value = 0.104256763104
$ws.Range("R4").Value = 0.9383108679359999
$ws.Range("S4").Value = 0.0001628258632370687
$ws.Range("T4").Value = 0.000178403460531506

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Pros1"
$ws.Range("C5").Value = "Tyro3"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 23.959728
$ws.Range("H5").Value = 71.879184
$ws.Range("I5").Value = 0.1471048789571275
$ws.Range("J5").Value = 0.1548696879699961
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4616195
$ws.Range("N5").Value = 0.923239
$ws.Range("O5").Value = 0.117424333800012
$ws.Range("P5").Value = 0.0814718136816918
$ws.Range("Q5").Value = 11.060277659496
$ws.Range("R5").Value = 66.361665956976
$ws.Range("S5").Value = 0.01727369241027211
$ws.Range("T5").Value = 0.01261751436323327

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Pros1"
$ws.Range("C6").Value = "Tyro3"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 85.25773366666668
$ws.Range("H6").Value = 255.773201
$ws.Range("I6").Value = 0.5234545480313474
$ws.Range("J6").Value = 0.5510846621458181
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.2063693333333333
$ws.Range("N6").Value = 0.619108
$ws.Range("O6").Value = 0.05249514260861875
$ws.Range("P6").Value = 0.05463357984752036
$ws.Range("Q6").Value = 17.59458165830089
$ws.Range("R6").Value = 158.351234924708
$ws.Range("S6").Value = 0.02747882114803566
$ws.Range("T6").Value = 0.03010772789208733

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Pros1"
$ws.Range("C7").Value = "Tyro3"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 85.25773366666668
$ws.Range("H7").Value = 255.773201
$ws.Range("I7").Value = 0.5234545480313474
$ws.Range("J7").Value = 0.5510846621458181
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.258868
$ws.Range("N7").Value = 9.776603999999999
$ws.Range("O7").Value = 0.8289736543672389
$ws.Range("P7").Value = 0.8627426479250582
$ws.Range("Q7").Value = 277.8436999988227
$ws.Range("R7").Value = 2500.593299989404
$ws.Range("S7").Value = 0.4339300295766974
$ws.Range("T7").Value = 0.4754442406505692

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Pros1"
$ws.Range("C8").Value = "Tyro3"
$ws.Range("D8").Value = "M1"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 85.25773366666668
$ws.Range("H8").Value = 255.773201
$ws.Range("I8").Value = 0.5234545480313474
$ws.Range("J8").Value = 0.5510846621458181
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.004351333333333333
$ws.Range("N8").Value = 0.013054
$ws.Range("O8").Value = 0.001106869224130377
$ws.Range("P8").Value = 0.001151958545729551
$ws.Range("Q8").Value = 0.3709848184282222
$ws.Range("R8").Value = 3.338863365854
$ws.Range("S8").Value = 0.0005793957294469745
$ws.Range("T8").Value = 0.0006348266859793574

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Pros1"
$ws.Range("C9").Value = "Tyro3"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 85.25773366666668
$ws.Range("H9").Value = 255.773201
$ws.Range("I9").Value = 0.5234545480313474
$ws.Range("J9").Value = 0.5510846621458181
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4616195
$ws.Range("N9").Value = 0.923239
$ws.Range("O9").Value = 0.117424333800012
$ws.Range("P9").Value = 0.0814718136816918
$ws.Range("Q9").Value = 39.35663238633984
$ws.Range("R9").Value = 236.139794318039
$ws.Range("S9").Value = 0.06146630157716738
$ws.Range("T9").Value = 0.04489786691718217

# Row 10
$ws.Range("A10").Value = "M1"
$ws.Range("B10").Value = "Pros1"
$ws.Range("C10").Value = "Tyro3"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 13.737404
$ws.Range("H10").Value = 41.212212
$ws.Range("I10").Value = 0.08434315917965177
$ws.Range("J10").Value = 0.088795142874651
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2063693333333333
$ws.Range("N10").Value = 0.619108
$ws.Range("O10").Value = 0.05249514260861875
$ws.Range("P10").Value = 0.05463357984752036
$ws.Range("Q10").Value = 2.834978905210666
$ws.Range("R10").Value = 25.514810146896
$ws.Range("S10").Value = 0.004427606169197251
$ws.Range("T10").Value = 0.004851196528314224

# Row 11
$ws.Range("A11").Value = "M1"
$ws.Range("B11").Value = "Pros1"
$ws.Range("C11").Value = "Tyro3"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 13.737404
$ws.Range("H11").Value = 41.212212
$ws.Range("I11").Value = 0.08434315917965177
$ws.Range("J11").Value = 0.088795142874651
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 3.258868
$ws.Range("N11").Value = 9.776603999999999
$ws.Range("O11").Value = 0.8289736543672389
$ws.Range("P11").Value = 0.8627426479250582
$ws.Range("Q11").Value = 44.76838629867199
$ws.Range("R11").Value = 402.915476688048
$ws.Range("S11").Value = 0.06991825688603366
$ws.Range("T11").Value = 0.07660735668656027

# Row 12
$ws.Range("A12").Value = "M1"
$ws.Range("B12").Value = "Pros1"
$ws.Range("C12").Value = "Tyro3"
$ws.Range("D12").Value = "M1"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 13.737404
$ws.Range("H12").Value = 41.212212
$ws.Range("I12").Value = 0.08434315917965177
$ws.Range("J12").Value = 0.088795142874651
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.004351333333333333
$ws.Range("N12").Value = 0.013054
$ws.Range("O12").Value = 0.001106869224130377
$ws.Range("P12").Value = 0.001151958545729551
$ws.Range("Q12").Value = 0.05977602393866666
$ws.Range("R12").Value = 0.537984215448
$ws.Range("S12").Value = 0.00009335684716188602
$ws.Range("T12").Value = 0.0001022883236537306

# Row 13
$ws.Range("A13").Value = "M1"
$ws.Range("B13").Value = "Pros1"
$ws.Range("C13").Value = "Tyro3"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 13.737404
$ws.Range("H13").Value = 41.212212
$ws.Range("I13").Value = 0.08434315917965177
$ws.Range("J13").Value = 0.088795142874651
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4616195
$ws.Range("N13").Value = 0.923239
$ws.Range("O13").Value = 0.117424333800012
$ws.Range("P13").Value = 0.0814718136816918
$ws.Range("Q13").Value = 6.341453565778
$ws.Range("R13").Value = 38.048721394668
$ws.Range("S13").Value = 0.00990393927725898
$ws.Range("T13").Value = 0.00723430133612277

# Row 14
$ws.Range("A14").Value = "M2"
$ws.Range("B14").Value = "Pros1"
$ws.Range("C14").Value = "Tyro3"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 15.421724
$ws.Range("H14").Value = 46.26517200000001
$ws.Range("I14").Value = 0.09468433207297799
$ws.Range("J14").Value = 0.09968216600119166
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.2063693333333333
$ws.Range("N14").Value = 0.619108
$ws.Range("O14").Value = 0.05249514260861875
$ws.Range("P14").Value = 0.05463357984752036
$ws.Range("Q14").Value = 3.182570900730667
$ws.Range("R14").Value = 28.643138106576
$ws.Range("S14").Value = 0.004970467514972795
$ws.Range("T14").Value = 0.005445993575599884

# Row 15
$ws.Range("A15").Value = "M2"
$ws.Range("B15").Value = "Pros1"
$ws.Range("C15").Value = "Tyro3"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 15.421724
$ws.Range("H15").Value = 46.26517200000001
$ws.Range("I15").Value = 0.09468433207297799
$ws.Range("J15").Value = 0.09968216600119166
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 3.258868
$ws.Range("N15").Value = 9.776603999999999
$ws.Range("O15").Value = 0.8289736543672389
$ws.Range("P15").Value = 0.8627426479250582
$ws.Range("Q15").Value = 50.257362848432
$ws.Range("R15").Value = 452.316265635888
$ws.Range("S15").Value = 0.07849081676985774
$ws.Range("T15").Value = 0.0860000558467733

# Row 16
$ws.Range("A16").Value = "M2"
$ws.Range("B16").Value = "Pros1"
$ws.Range("C16").Value = "Tyro3"
$ws.Range("D16").Value = "M1"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 15.421724
$ws.Range("H16").Value = 46.26517200000001
$ws.Range("I16").Value = 0.09468433207297799
$ws.Range("J16").Value = 0.09968216600119166
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.004351333333333333
$ws.Range("N16").Value = 0.013054
$ws.Range("O16").Value = 0.001106869224130377
$ws.Range("P16").Value = 0.001151958545729551
$ws.Range("Q16").Value = 0.06710506169866667
$ws.Range("R16").Value = 0.6039455552880001
$ws.Range("S16").Value = 0.0001048031731789201
$ws.Range("T16").Value = 0.0001148297229819044

# Row 17
$ws.Range("A17").Value = "M2"
$ws.Range("B17").Value = "Pros1"
$ws.Range("C17").Value = "Tyro3"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 15.421724
$ws.Range("H17").Value = 46.26517200000001
$ws.Range("I17").Value = 0.09468433207297799
$ws.Range("J17").Value = 0.09968216600119166
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4616195
$ws.Range("N17").Value = 0.923239
$ws.Range("O17").Value = 0.117424333800012
$ws.Range("P17").Value = 0.0814718136816918
$ws.Range("Q17").Value = 7.118968522018002
$ws.Range("R17").Value = 42.71381113210801
$ws.Range("S17").Value = 0.01111824461496855
$ws.Range("T17").Value = 0.00812128685583656

# Row 18
$ws.Range("A18").Value = "sCs"
$ws.Range("B18").Value = "Pros1"
$ws.Range("C18").Value = "Tyro3"
$ws.Range("D18").Value = "ECs"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 24.498552
$ws.Range("H18").Value = 48.997104
$ws.Range("I18").Value = 0.1504130817588954
$ws.Range("J18").Value = 0.1055683410083432
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 0.2063693333333333
$ws.Range("N18").Value = 0.619108
$ws.Range("O18").Value = 0.05249514260861875
$ws.Range("P18").Value = 0.05463357984752036
$ws.Range("Q18").Value = 5.055749843871999
$ws.Range("R18").Value = 30.334499063232
$ws.Range("S18").Value = 0.007895956177135045
$ws.Range("T18").Value = 0.005767576387849576

# Row 19
$ws.Range("A19").Value = "sCs"
$ws.Range("B19").Value = "Pros1"
$ws.Range("C19").Value = "Tyro3"
$ws.Range("D19").Value = "FAPs"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 24.498552
$ws.Range("H19").Value = 48.997104
$ws.Range("I19").Value = 0.1504130817588954
$ws.Range("J19").Value = 0.1055683410083432
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 3.258868
$ws.Range("N19").Value = 9.776603999999999
$ws.Range("O19").Value = 0.8289736543672389
$ws.Range("P19").Value = 0.8627426479250582
$ws.Range("Q19").Value = 79.83754715913599
$ws.Range("R19").Value = 479.025282954816
$ws.Range("S19").Value = 0.1246884820503098
$ws.Range("T19").Value = 0.09107831005859351

# Row 20
$ws.Range("A20").Value = "sCs"
$ws.Range("B20").Value = "Pros1"
$ws.Range("C20").Value = "Tyro3"
$ws.Range("D20").Value = "M1"
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 24.498552
$ws.Range("H20").Value = 48.997104
$ws.Range("I20").Value = 0.1504130817588954
$ws.Range("J20").Value = 0.1055683410083432
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.004351333333333333
$ws.Range("N20").Value = 0.013054
$ws.Range("O20").Value = 0.001106869224130377
$ws.Range("P20").Value = 0.001151958545729551
$ws.Range("Q20").Value = 0.106601365936
$ws.Range("R20").Value = 0.639608195616
$ws.Range("S20").Value = 0.0001664876111055274
$ws.Range("T20").Value = 0.0001216103525830523

# Row 21
$ws.Range("A21").Value = "sCs"
$ws.Range("B21").Value = "Pros1"
$ws.Range("C21").Value = "Tyro3"
$ws.Range("D21").Value = "sCs"
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 24.498552
$ws.Range("H21").Value = 48.997104
$ws.Range("I21").Value = 0.1504130817588954
$ws.Range("J21").Value = 0.1055683410083432
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.4616195
$ws.Range("N21").Value = 0.923239
$ws.Range("O21").Value = 0.117424333800012
$ws.Range("P21").Value = 0.0814718136816918
$ws.Range("Q21").Value = 11.309009324964
$ws.Range("R21").Value = 45.236037299856
$ws.Range("S21").Value = 0.01766215592034503
$ws.Range("T21").Value = 0.008600844209317041
